$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx tracking numbers for column P (rows 2-26). These are long
# numeric-looking strings that must stay stored as TEXT (matching the
# existing shared-string cells), not get auto-converted to a number by
# Excel's input parser. We briefly force text entry via a "@" number
# format, then reset the cell style back to "Normal" so the cell keeps
# its original (unstyled) appearance, matching the rest of the column.
$trackingNumbers = @{
    2  = "320017963792"
    3  = "320017963807"
    4  = "320017963830"
    5  = "320017963851"
    6  = "320017963895"
    7  = "320017963910"
    8  = "320017963943"
    9  = "320017963965"
    10 = "320017963998"
    11 = "320017964012"
    12 = "320017964056"
    13 = "320017964078"
    14 = "320017964104"
    15 = "320017964126"
    16 = "320017964159"
    17 = "320017964170"
    18 = "320017964218"
    19 = "320017964230"
    20 = "320017964284"
    21 = "320017964300"
    22 = "320017964332"
    23 = "320017964343"
    24 = "320017964354"
    25 = "320017964365"
    26 = "320017964376"
}

foreach ($row in $trackingNumbers.Keys) {
    $cell = $ws.Range("P$row")
    $cell.NumberFormat = "@"
    $cell.Value = $trackingNumbers[$row]
    $cell.Style = "Normal"
}

# Row 24 also changes its ActualRate (Q) and Result (R) values: the
# rerun shipment's actual charge no longer matches the expected rate,
# so the result flips from PASS to FAIL.
$q24 = $ws.Range("Q24")
$q24.NumberFormat = "@"
$q24.Value = "$253.80"
$q24.Style = "Normal"

$ws.Range("R24").Value = "FAIL"
